{"js": "// Update the date header and the 25 division problems in the table,\n// matching the commit's output (generated at 9a8706d).\nconst replacements = [\n  [\"2024-01-11 Thursday\", \"2024-01-12 Friday\"],\n  [\"640\u00f76=\", \"373\u00f75=\"],\n  [\"300\u00f72=\", \"995\u00f79=\"],\n  [\"419\u00f76=\", \"276\u00f74=\"],\n  [\"620\u00f78=\", \"176\u00f73=\"],\n  [\"496\u00f76=\", \"408\u00f72=\"],\n  [\"855\u00f77=\", \"116\u00f72=\"],\n  [\"195\u00f76=\", \"869\u00f77=\"],\n  [\"961\u00f73=\", \"441\u00f72=\"],\n  [\"915\u00f79=\", \"631\u00f75=\"],\n  [\"626\u00f75=\", \"619\u00f74=\"],\n  [\"472\u00f79=\", \"826\u00f77=\"],\n  [\"719\u00f76=\", \"164\u00f78=\"],\n  [\"894\u00f79=\", \"453\u00f73=\"],\n  [\"938\u00f73=\", \"985\u00f76=\"],\n  [\"757\u00f74=\", \"619\u00f74=\"],\n  [\"606\u00f76=\", \"315\u00f77=\"],\n  [\"579\u00f75=\", \"163\u00f74=\"],\n  [\"215\u00f77=\", \"768\u00f78=\"],\n  [\"129\u00f73=\", \"748\u00f77=\"],\n  [\"342\u00f72=\", \"124\u00f77=\"],\n  [\"107\u00f77=\", \"146\u00f78=\"],\n  [\"465\u00f76=\", \"720\u00f77=\"],\n  [\"634\u00f74=\", \"953\u00f74=\"],\n  [\"221\u00f77=\", \"666\u00f79=\"],\n  [\"555\u00f74=\", \"567\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and the 25 division problems in the table,\n# matching the commit's output (generated at 9a8706d).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-11 Thursday\", \"2024-01-12 Friday\"),\n    @(\"640\u00f76=\", \"373\u00f75=\"),\n    @(\"300\u00f72=\", \"995\u00f79=\"),\n    @(\"419\u00f76=\", \"276\u00f74=\"),\n    @(\"620\u00f78=\", \"176\u00f73=\"),\n    @(\"496\u00f76=\", \"408\u00f72=\"),\n    @(\"855\u00f77=\", \"116\u00f72=\"),\n    @(\"195\u00f76=\", \"869\u00f77=\"),\n    @(\"961\u00f73=\", \"441\u00f72=\"),\n    @(\"915\u00f79=\", \"631\u00f75=\"),\n    @(\"626\u00f75=\", \"619\u00f74=\"),\n    @(\"472\u00f79=\", \"826\u00f77=\"),\n    @(\"719\u00f76=\", \"164\u00f78=\"),\n    @(\"894\u00f79=\", \"453\u00f73=\"),\n    @(\"938\u00f73=\", \"985\u00f76=\"),\n    @(\"757\u00f74=\", \"619\u00f74=\"),\n    @(\"606\u00f76=\", \"315\u00f77=\"),\n    @(\"579\u00f75=\", \"163\u00f74=\"),\n    @(\"215\u00f77=\", \"768\u00f78=\"),\n    @(\"129\u00f73=\", \"748\u00f77=\"),\n    @(\"342\u00f72=\", \"124\u00f77=\"),\n    @(\"107\u00f77=\", \"146\u00f78=\"),\n    @(\"465\u00f76=\", \"720\u00f77=\"),\n    @(\"634\u00f74=\", \"953\u00f74=\"),\n    @(\"221\u00f77=\", \"666\u00f79=\"),\n    @(\"555\u00f74=\", \"567\u00f72=\")\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n# MatchCase:=True keeps each swap targeted at its exact source string;\n# Replace:=2 is wdReplaceAll (harmless here since every string is unique).\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
